# lcms-metadata.xlsx: add an "is_targeted list" sheet (TRUE/FALSE) and point
# the "is_targeted" column's data validation at it instead of the old
# hard-coded "TRUE,FALSE" literal list.

$wb = $excel.ActiveWorkbook

$mainSheet = $wb.Worksheets.Item("Export as TSV")
$afterSheet = $wb.Worksheets.Item("analyte_class list")

# Insert the new lookup sheet right after "analyte_class list" (so it lands
# before "ms_source list", matching the rest of the "<field> list" sheets).
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "is_targeted list"

# Write literal TRUE / FALSE as TEXT (not boolean) values. Assigning the
# string directly via .Value auto-coerces to a real boolean cell, so instead
# write a text formula and flatten it back down to a static value via
# copy / paste-special-values - this keeps the cells as plain shared-string
# text cells (t="s"), matching the other lookup-list sheets.
$newSheet.Range("A1").Formula = '="TRUE"'
$newSheet.Range("A2").Formula = '="FALSE"'
$newSheet.Range("A1:A2").Copy()
$newSheet.Range("A1:A2").PasteSpecial(-4163)

# Point column N's validation at the new list instead of the inline
# "TRUE,FALSE" literal, and refresh the title/message to match the other
# list-backed columns.
$boolRange = $mainSheet.Range("N2:N1048576")
$boolRange.Validation.Formula1 = "='is_targeted list'!`$A`$1:`$A`$2"
$boolRange.Validation.ErrorTitle = "Value must come from list"
$boolRange.Validation.ErrorMessage = "Value must be one of: TRUE / FALSE."

# Leave the workbook focused back on the main data-entry sheet.
$mainSheet.Select()
